# This workbook's source feed re-ordered its records: the observation that
# used to be the last data row (row 7 - "Spindelblomster" / Neottia cordata,
# reported 2022-08-14 by Ola Löfgren at "Nöjden, 2,6 km V från befintligt
# hus, Jmt") is now the first data row (row 2), and every other record
# shifts down by one row (old row 2 -> new row 3, old 3 -> 4, old 4 -> 5,
# old 5 -> 6, old 6 -> 7). Apply that re-shuffle by writing each row's
# values into its new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2 (previously row 7: Spindelblomster / Neottia cordata, Nöjden) ---
$ws.Range("A2").Value = 103672936
$ws.Range("B2").Value = 96354
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 221952
$ws.Range("F2").Value = "Spindelblomster"
$ws.Range("G2").Value = "Neottia cordata"
$ws.Range("H2").Value = "(L.) Rich."
$ws.Range("P2").Value = "Nöjden, 2,6 km V från befintligt hus, Jmt"
$ws.Range("Q2").Value = 444643.6316846627
$ws.Range("R2").Value = 7094865.722657905
# leading apostrophe keeps these as literal text (not auto-converted to a date serial)
$ws.Range("Y2").Value = "'2022-08-14"
$ws.Range("AA2").Value = "'2022-08-14"
$ws.Range("AI2").Value = "Vid bäck i granskog"
$ws.Range("AW2").Value = "Ola Löfgren"
$ws.Range("AX2").Value = "Ola Löfgren"
# cells that no longer carry (even empty) content in the new row 2
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("AF2").ClearContents()

# --- New row 3 (previously row 2: Gränsticka / Phellopilus nigrolimitatus) ---
$ws.Range("A3").Value = 101362742
$ws.Range("B3").Value = 89406
$ws.Range("E3").Value = 1204
$ws.Range("F3").Value = "Gränsticka"
$ws.Range("G3").Value = "Phellopilus nigrolimitatus"
$ws.Range("H3").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."

# --- New row 4 (previously row 3: Gammelgransskål / Pseudographis pinicola) ---
$ws.Range("A4").Value = 101362755
$ws.Range("B4").Value = 81236
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1312
$ws.Range("F4").Value = "Gammelgransskål"
$ws.Range("G4").Value = "Pseudographis pinicola"
$ws.Range("H4").Value = "(Nyl.) Rehm"

# --- New row 5 (previously row 4: Spindelblomster / Neottia cordata, Nöjdfjället) ---
$ws.Range("A5").Value = 101362747
$ws.Range("B5").Value = 96354
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221952
$ws.Range("F5").Value = "Spindelblomster"
$ws.Range("G5").Value = "Neottia cordata"
$ws.Range("H5").Value = "(L.) Rich."
$ws.Range("L5").ClearContents()
$ws.Range("Q5").Value = 444383.9599903998
$ws.Range("R5").Value = 7094782.483513337

# --- New row 6 (previously row 5: Garnlav / Alectoria sarmentosa) ---
$ws.Range("A6").Value = 101362767
$ws.Range("B6").Value = 77506
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 444364.4024028718
$ws.Range("R6").Value = 7094746.371122397

# --- New row 7 (previously row 6: Doftskinn / Cystostereum murrayi) ---
$ws.Range("A7").Value = 101362741
$ws.Range("B7").Value = 85703
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 510
$ws.Range("F7").Value = "Doftskinn"
$ws.Range("G7").Value = "Cystostereum murrayi"
$ws.Range("H7").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("P7").Value = "Nöjdfjället, Jmt"
$ws.Range("Q7").Value = 444383.9599903998
$ws.Range("R7").Value = 7094782.483513337
$ws.Range("Y7").Value = "'2022-06-02"
$ws.Range("AA7").Value = "'2022-06-02"
$ws.Range("AW7").Value = "Andreas Öster"
$ws.Range("AX7").Value = "Andreas Öster"
# cells that carried (empty) content before but no longer apply to this row
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AI7").ClearContents()
